{"js": "// The trailing paragraph of the document is gibberish test text\n// (\"djlsajdl...\") that was accidentally left in the manuscript. The\n// author removed it and merged what's left (an otherwise-empty\n// paragraph) back into the previous paragraph, which ends in\n// \"...sent\u00eda como el aire llega con dificultad\". The \"_GoBack\" bookmark\n// that lived in the gibberish paragraph has to survive the merge, now\n// sitting at the end of the \"dificultad\" paragraph.\n\nconst body = context.document.body;\n\n// 1) Remove the gibberish run itself.\nconst gibberish = body.search(\n  \"djlsajdlajdlasjdlajdlajdllkasjdlkasjdlkasjsdlaskjdlakjdlaslkdjaslkdjalskdjaslkdjas\",\n  { matchCase: true }\n);\nawait context.sync();\n\nif (gibberish.items.length > 0) {\n  gibberish.items[0].delete();\n  await context.sync();\n}\n\n// 2) Re-create the \"_GoBack\" bookmark right after \"dificultad\" (i.e. at\n//    the end of what will become the merged paragraph) before we drop\n//    the paragraph that currently owns it.\nconst anchor = body.search(\"dificultad\", { matchCase: true });\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const afterAnchor = anchor.items[0].getRange(\"After\");\n  afterAnchor.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 3) Drop the now-empty trailing paragraph (and with it, its old\n//    bookmark/proofErr marks), effectively merging it into the\n//    previous paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nif (count > 0) {\n  const lastParagraph = paragraphs.items[count - 1];\n  lastParagraph.delete();\n  await context.sync();\n}\n", "ps1": "# The trailing paragraph of the document is gibberish test text\n# (\"djlsajdl...\") that was accidentally left in the manuscript. The\n# author removed it and merged what's left (an otherwise-empty\n# paragraph) back into the previous paragraph, which ends in\n# \"...sent\u00eda como el aire llega con dificultad\". The \"_GoBack\" bookmark\n# that lived in the gibberish paragraph has to survive the merge, now\n# sitting at the end of the \"dificultad\" paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the gibberish run itself.\n$gibberishText = \"djlsajdlajdlasjdlajdlajdllkasjdlkasjdlkasjsdlaskjdlakjdlaslkdjaslkdjalskdjaslkdjas\"\n$delRange = $d.Content\n$foundGibberish = $delRange.Find.Execute($gibberishText)\nif ($foundGibberish) {\n    $delRange.Delete()\n}\n\n# 2) Re-create the \"_GoBack\" bookmark right after \"dificultad\" (i.e. at\n#    the end of what will become the merged paragraph) before we drop\n#    the paragraph that currently owns it. A temporary placeholder\n#    character is used so the bookmark can be anchored precisely, then\n#    the placeholder is removed, leaving a zero-length bookmark.\n$anchorRange = $d.Content\n$foundAnchor = $anchorRange.Find.Execute(\"dificultad\")\nif ($foundAnchor) {\n    $anchorRange.Collapse([Microsoft.Office.Interop.Word.WdCollapseDirection]::wdCollapseEnd)\n    $anchorRange.InsertAfter(\"x\")\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks(\"_GoBack\").Delete()\n    }\n    $d.Bookmarks.Add(\"_GoBack\", $anchorRange)\n    $anchorRange.Text = \"\"\n}\n\n# 3) Drop the now-empty trailing paragraph (and with it, its old\n#    bookmark/proofErr marks), effectively merging it into the\n#    previous paragraph.\n$count = $d.Paragraphs.Count\nif ($count -gt 0) {\n    $d.Paragraphs.Item($count).Range.Delete()\n}\n"}
